$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Synonyms" label (row 5, column D) to "NETs - General"
$ws.Range("D5").Value = "NETs - General"

# Renumber the trailing rank/index values in column A
$ws.Range("A7").Value = 6
$ws.Range("A9").Value = 7
$ws.Range("A10").Value = 8
